{"js": "// Insert the 3 new \"Move the ... first\" sub-bullets right after the\n// \"ID potential solutions\" bullet (the first one, under the first \"Break\n// the problem apart\" section), matching the commit:\n// \"I've added the 3 potential solutions to this equation.\"\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\n// Locate the first paragraph whose text is exactly \"ID potential solutions\".\nlet anchor = null;\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  if (paragraphs.items[i].text === \"ID potential solutions\") {\n    anchor = paragraphs.items[i];\n    break;\n  }\n}\n\nif (!anchor) {\n  throw new Error('Could not find the \"ID potential solutions\" paragraph.');\n}\n\nconst newLines = [\"Move the cat first\", \"Move the seed first\", \"Move the bird first\"];\n\n// Insert each new paragraph right after the anchor, in order, so they end\n// up appearing \"Move the cat first\", \"Move the seed first\", \"Move the bird\n// first\" directly beneath \"ID potential solutions\".\nlet insertAfter = anchor.getRange(\"End\");\nconst newParaRanges = [];\nfor (const line of newLines) {\n  const newRange = insertAfter.insertParagraph(line, \"After\");\n  newParaRanges.push(newRange);\n  insertAfter = newRange;\n}\nawait context.sync();\n\n// Promote each newly inserted paragraph to the second list level (ilvl=1),\n// matching its sibling bullets (\"The Cat\" / \"The Bird\" / \"The Seed\").\nfor (const r of newParaRanges) {\n  const p = r.paragraphs.getFirst();\n  p.listItem.level = 1;\n}\nawait context.sync();\n", "ps1": "# Insert the 3 new \"Move the ... first\" sub-bullets right after the\n# \"ID potential solutions\" bullet (the first one, under the first \"Break\n# the problem apart\" section), matching the commit:\n# \"I've added the 3 potential solutions to this equation.\"\n\n$d = $word.ActiveDocument\n\n# Locate the first paragraph whose text is exactly \"ID potential solutions\".\n$anchorIndex = -1\nfor ($i = 1; $i -le $d.Paragraphs.Count; $i++) {\n    $t = $d.Paragraphs.Item($i).Range.Text.TrimEnd([char]13, [char]7)\n    if ($t -eq \"ID potential solutions\") {\n        $anchorIndex = $i\n        break\n    }\n}\n\nif ($anchorIndex -eq -1) {\n    throw \"Could not find the 'ID potential solutions' paragraph.\"\n}\n\n$lines = @(\"Move the cat first\", \"Move the seed first\", \"Move the bird first\")\n\n# Insert each new paragraph right after the anchor, in order, so they end up\n# appearing \"Move the cat first\", \"Move the seed first\", \"Move the bird\n# first\" directly beneath \"ID potential solutions\". Each new paragraph is\n# promoted to the second list level (ilvl=1 / ListLevelNumber=2), matching\n# its sibling bullets (\"The Cat\" / \"The Bird\" / \"The Seed\").\n$insertIndex = $anchorIndex\nforeach ($line in $lines) {\n    $p = $d.Paragraphs.Item($insertIndex)\n    $p.Range.InsertParagraphAfter()\n    $insertIndex = $insertIndex + 1\n    $newPara = $d.Paragraphs.Item($insertIndex)\n    $newPara.Range.Text = $line\n    $newPara.Range.ListFormat.ListLevelNumber = 2\n}\n"}
